# Add a new task row (row 51) to the TodoList sheet, describing the
# "DlogGroup-simultaneous exponentiations" bug/task, right after the last
# existing row (row 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Start from the formatting of the row directly above (row 49 uses the same
# banded row style that row 51 should end up with) so the new row matches
# the alternating green/orange row banding used throughout the sheet.
$ws.Range("A49:G49").Copy()
$ws.Range("A51:G51").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new task's data.
$ws.Range("A51").Value = "DlogGroup-simultaneous exponentiations"
$ws.Range("B51").Value = "Bug Fix"
$ws.Range("C51").Value = "This is not really a bug. At the moment we have implemented this function in a naïve way, without the optimization suggested in the book of applied crypto. See if we can optimize."
$ws.Range("D51").Value = "2/20/2012"
$ws.Range("E51").Value = "Yael"
$ws.Range("F51").Value = "High"
$ws.Range("G51").Value = ""

# The description wraps onto two lines at this column width, so the row is
# taller than the single-line rows.
$ws.Rows.Item(51).RowHeight = 30

# Leave the newly added row selected, like a user would after typing the
# last entry and pressing Enter on the whole row.
[void]$ws.Rows.Item(51).EntireRow.Select()
